$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row (row 1) with new columns P and Q (copy O1's format: bold, bordered, centered)
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# For each data row (2..25), fix columns I, K, M, O (swap 1/2) and add P, Q (=2)
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O: 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P: new = 2
    $ws.Cells.Item($r, 17).Value = 2   # Q: new = 2
}
